# Daily attendance processing - 2025-11-05 10:25:09
# Reverse the order of the comma-separated "Recorded By" names in column G
# for every data row on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        $n = $parts.Count
        $reversedParts = $parts[($n - 1)..0]
        $newVal = [string]::Join(", ", $reversedParts)
        $cell.Value2 = $newVal
    }
}
